# feat: update foreach/endrow/endloop with new behaviour
#
# Adds a new "#! END_ROW true" marker to a new column J on both sheets,
# used next to FOR_EACH / CONTINUE rows (loop rows that repeat), while the
# existing "#! END_ROW" marker is kept for the END_LOOP row (last row of
# the block).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet1: single-iteration FOR_EACH / END_LOOP block
$ws1.Range("J2").Value = "#! END_ROW true"
$ws1.Range("J3").Value = "#! END_ROW"

# Sheet2: FOR_EACH / CONTINUE / END_LOOP block
$ws2.Range("J2").Value = "#! END_ROW true"
$ws2.Range("J3").Value = "#! END_ROW true"
$ws2.Range("J4").Value = "#! END_ROW"

# Restore the view state: Sheet1 active with selection at G12,
# Sheet2 inactive with selection at I2.
$ws2.Activate()
$ws2.Range("I2").Select()
$ws1.Activate()
$ws1.Range("G12").Select()
